# Insert a new data row at row 168 (pushing the existing rows 168-185 down
# to 169-186) and populate the new row with its values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 168.. down by one row.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new record's data.
$ws.Range("A168").Value = 11
$ws.Range("B168").Value = "Vega Monumental Concepción"
$ws.Range("C168").Value = "Bíobío"
$ws.Range("D168").Value = 44946
$ws.Range("E168").Value = 8
$ws.Range("F168").Value = 100112032
$ws.Range("G168").Value = "Zapallo italiano"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 350
$ws.Range("K168").Value = 4000
$ws.Range("L168").Value = 5000
$ws.Range("M168").Value = 4429
$ws.Range("N168").Value = "`$/caja 50 unidades"
$ws.Range("O168").Value = "Región de O'Higgins"
$ws.Range("P168").Value = 89
$ws.Range("Q168").Value = 50
$ws.Range("R168").Value = "Hortaliza"
